$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row before row 20, pushing the old rows 20-22 down to 21-23.
$ws.Rows("20:20").Insert()

# Copy the (now shifted) row 21 formatting onto the new, blank row 20 so the
# new row picks up the same shaded style (fill) that row used to have.
$ws.Range("A21:G21").Copy()
$ws.Range("A20:G20").PasteSpecial(-4122)

# Fill in the values for the newly inserted row 20.
$ws.Range("A20").Value = "YES"
$ws.Range("B20").Value = "server12"
$ws.Range("C20").Value = "Availability"
$ws.Range("D20").Value = "CRITICAL"
$ws.Range("E20").Value = "<"
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = "Availability"

# Update the view: scroll a couple of rows further and select the full row 20.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 6
$ws.Range("A20:XFD20").Select()
